$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("C2").Value = 12
$ws.Range("C3").Value = 11

# Update the current selection to match the diff
$ws.Range("B2:C3").Select()
